# Generate Report for Handback
# Update the "Correspond Handoff Datetime" and "Correspond Handback DateTime"
# values for the first data row (the 8e2b4336... file) on both the zh-cn and
# de-de sheets, reflecting a newer handback report run. The second data row
# (e0390748...) keeps its original timestamps.

$wb = $excel.ActiveWorkbook

$ws_zhcn = $wb.Worksheets.Item("zh-cn")
$ws_zhcn.Range("E2").Value = "2016-03-14 04:30:56"
$ws_zhcn.Range("H2").Value = "2016-03-14 04:31:13"

$ws_dede = $wb.Worksheets.Item("de-de")
$ws_dede.Range("E2").Value = "2016-03-14 04:30:59"
$ws_dede.Range("H2").Value = "2016-03-14 04:31:18"
